$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from an existing header cell (H1) so bold/border/alignment match
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for columns I and J, rows 2-13
$values = @{
    2  = @(3, 7)
    3  = @(1, 5)
    4  = @(2, 6)
    5  = @(1, 5)
    6  = @(1, 6)
    7  = @(8, 9)
    8  = @(3, 8)
    9  = @(3, 6)
    10 = @(5, 8)
    11 = @(7, 9)
    12 = @(7, 7)
    13 = @(7, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
